$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 459.0625
$ws.Range("I41").Value = 218.2
$ws.Range("J41").Value = 860.5
$ws.Range("K41").Value = 218.2
$ws.Range("L41").Value = 860.5
$ws.Range("M41").Value = 221.8
$ws.Range("N41").Value = -1740.5
$ws.Range("H62").Value = 2706.6
$ws.Range("I62").Value = 2634
$ws.Range("K62").Value = 2634
$ws.Range("M62").Value = -2010
$ws.Range("H65").Value = 2706.6
$ws.Range("I65").Value = 2634
$ws.Range("K65").Value = 13170
$ws.Range("M65").Value = -10050
$ws.Range("H70").Value = 9464
$ws.Range("I70").Value = 5850
$ws.Range("J70").Value = 10576
$ws.Range("K70").Value = 17550
$ws.Range("L70").Value = 31728
$ws.Range("M70").Value = -17280
$ws.Range("N70").Value = -32268
$ws.Range("H73").Value = 9464
$ws.Range("I73").Value = 5850
$ws.Range("J73").Value = 10576
$ws.Range("K73").Value = 17550
$ws.Range("L73").Value = 31728
$ws.Range("M73").Value = -16614
$ws.Range("N73").Value = -33600
$ws.Range("H74").Value = 6197.143
$ws.Range("I74").Value = 4707.5
$ws.Range("K74").Value = 4707.5
$ws.Range("M74").Value = -3771.5
$ws.Range("H77").Value = 6197.143
$ws.Range("I77").Value = 4707.5
$ws.Range("K77").Value = 23537.5
$ws.Range("M77").Value = -18857.5
$ws.Range("H99").Value = 1594.2222
$ws.Range("I99").Value = 241.66667
$ws.Range("J99").Value = 4299.3335
$ws.Range("K99").Value = 725.00001
$ws.Range("L99").Value = 12898.0005
$ws.Range("M99").Value = 772.99999
$ws.Range("N99").Value = -15894.0005
$ws.Range("H101").Value = 2972.5264
$ws.Range("I101").Value = 510.4
$ws.Range("J101").Value = 5708.222
$ws.Range("K101").Value = 1531.2
$ws.Range("L101").Value = 17124.666
$ws.Range("M101").Value = 90.80000000000018
$ws.Range("N101").Value = -20368.666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7605.5625
$ws.Range("I61").Value = 6446.0967
$ws.Range("K61").Value = 6446.0967
$ws.Range("M61").Value = -6234.0967
$ws.Range("H74").Value = 4226.9287
$ws.Range("I74").Value = 3733.6667
$ws.Range("K74").Value = 3733.6667
$ws.Range("M74").Value = -2859.6667
$ws.Range("H77").Value = 4226.9287
$ws.Range("I77").Value = 3733.6667
$ws.Range("K77").Value = 18668.3335
$ws.Range("M77").Value = -14300.3335
$ws.Range("H102").Value = 1354.0667
$ws.Range("I102").Value = 1408.6154
$ws.Range("K102").Value = 1408.6154
$ws.Range("M102").Value = 213.3846000000001
$ws.Range("H110").Value = 3631.2856
$ws.Range("I110").Value = 3640
$ws.Range("K110").Value = 3640
$ws.Range("M110").Value = -1595
$ws.Range("H122").Value = 2606.923
$ws.Range("I122").Value = 1471
$ws.Range("J122").Value = 3316.875
$ws.Range("K122").Value = 4413
$ws.Range("L122").Value = 9950.625
$ws.Range("M122").Value = -1963
$ws.Range("N122").Value = -14850.625
$ws.Range("H132").Value = 3391.353
$ws.Range("I132").Value = 3010.1365
$ws.Range("K132").Value = 9030.4095
$ws.Range("M132").Value = -6500.4095
$ws.Range("H136").Value = 7605.5625
$ws.Range("I136").Value = 6446.0967
$ws.Range("K136").Value = 19338.2901
$ws.Range("M136").Value = -16788.2901

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2238
$ws.Range("I105").Value = 2150.6667
$ws.Range("K105").Value = 2150.6667
$ws.Range("M105").Value = -403.6667000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19499.7
$ws.Range("J41").Value = 35000
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35856
$ws.Range("H58").Value = 3469.1614
$ws.Range("I58").Value = 2744.348
$ws.Range("J58").Value = 5553
$ws.Range("K58").Value = 2744.348
$ws.Range("L58").Value = 5553
$ws.Range("M58").Value = -2541.348
$ws.Range("N58").Value = -5959
$ws.Range("H122").Value = 4100.2
$ws.Range("I122").Value = 4165.3687
$ws.Range("J122").Value = 3893.8333
$ws.Range("K122").Value = 12496.1061
$ws.Range("L122").Value = 11681.4999
$ws.Range("M122").Value = -10046.1061
$ws.Range("N122").Value = -16581.4999
$ws.Range("H134").Value = 4755.121
$ws.Range("I134").Value = 4257.5063
$ws.Range("J134").Value = 6720.7
$ws.Range("K134").Value = 12772.5189
$ws.Range("L134").Value = 20162.1
$ws.Range("M134").Value = -10237.5189
$ws.Range("N134").Value = -25232.1
$ws.Range("H136").Value = 3469.1614
$ws.Range("I136").Value = 2744.348
$ws.Range("J136").Value = 5553
$ws.Range("K136").Value = 8233.044
$ws.Range("L136").Value = 16659
$ws.Range("M136").Value = -5683.044
$ws.Range("N136").Value = -21759

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31654480
$ws.Range("I4").Value = 34894148
$ws.Range("J4").Value = 337666.66
$ws.Range("K4").Value = 104682444
$ws.Range("L4").Value = 1012999.98
$ws.Range("M4").Value = -104682332
$ws.Range("N4").Value = -1013223.98
$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685
$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908
$ws.Range("H117").Value = 2044.0555
$ws.Range("I117").Value = 1076.1666
$ws.Range("K117").Value = 3228.4998
$ws.Range("M117").Value = 213.5001999999999
$ws.Range("H129").Value = 2668.1
$ws.Range("I129").Value = 1041.5
$ws.Range("J129").Value = 5108
$ws.Range("K129").Value = 3124.5
$ws.Range("L129").Value = 15324
$ws.Range("M129").Value = 1875.5
$ws.Range("N129").Value = -25324

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 352.1111
$ws.Range("J107").Value = 444
$ws.Range("L107").Value = 444
$ws.Range("N107").Value = -4284
$ws.Range("H126").Value = 4637.278
$ws.Range("I126").Value = 2477.75
$ws.Range("K126").Value = 7433.25
$ws.Range("M126").Value = -4963.25
$ws.Range("H132").Value = 5003.8887
$ws.Range("I132").Value = 4445.1333
$ws.Range("K132").Value = 13335.3999
$ws.Range("M132").Value = -10805.3999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2059.3333
$ws.Range("J46").Value = 2408.182
$ws.Range("L46").Value = 2408.182
$ws.Range("N46").Value = -2784.182
$ws.Range("H61").Value = 1736.6923
$ws.Range("I61").Value = 1736.6923
$ws.Range("K61").Value = 1736.6923
$ws.Range("M61").Value = -1534.6923
$ws.Range("H93").Value = 1576.091
$ws.Range("I93").Value = 1498
$ws.Range("K93").Value = 1498
$ws.Range("M93").Value = -250
$ws.Range("H113").Value = 1736.6923
$ws.Range("I113").Value = 1736.6923
$ws.Range("K113").Value = 1736.6923
$ws.Range("M113").Value = 433.3077000000001
$ws.Range("H122").Value = 5755.8276
$ws.Range("I122").Value = 3716.25
$ws.Range("J122").Value = 6532.8096
$ws.Range("K122").Value = 11148.75
$ws.Range("L122").Value = 19598.4288
$ws.Range("M122").Value = -8698.75
$ws.Range("N122").Value = -24498.4288

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 655.26086
$ws.Range("I107").Value = 628.05554
$ws.Range("K107").Value = 1884.16662
$ws.Range("M107").Value = 35.83338000000003
$ws.Range("H122").Value = 298973.94
$ws.Range("I122").Value = 503806.06
$ws.Range("K122").Value = 1511418.18
$ws.Range("M122").Value = -1508968.18
